$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the six SMS-code labels (row 7 "Εκκλησία" is being retired).
$ws.Range("A1").Value = "Ιατρόν - Φαρμακεία"
$ws.Range("A2").Value = "Αγορά τροφίμων"
$ws.Range("A3").Value = "Τράπεζα Χρημάτων"
$ws.Range("A4").Value = "Βοήθεια κατ' οίκον"
$ws.Range("A5").Value = "Τελετή - Τέκνα στο γυμνάσιο"
$ws.Range("A6").Value = "Ασκηση - Περίπατος με οικόσιτο"

# A3/A5/A6 previously carried a distinct (unused) font/style; align them
# with A2's formatting so they share the same style record.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)

# Drop the now-removed 7th entry ("Εκκλησία").
$ws.Rows.Item(7).Delete()
